# Auto update Excel log
# Appends new sensor-log rows to the "PIR" and "Humidity" sheets.
#
# Note: columns that hold digit-led text which Excel's input-parser would
# otherwise "helpfully" reinterpret (the "yyyy-mm-dd" date strings, and the
# "NN.N%" humidity readings) are written as a literal `="..."` text formula
# and then flattened to a plain value via Copy / PasteSpecial(xlPasteValues)
# so the stored cell stays a plain string rather than becoming a date serial
# or a numeric percentage.

$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: append rows 167-180
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")

$pirData = @(
    @("18:35:08", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:09", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:11", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:16", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:21", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:26", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:31", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:36", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:41", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:46", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:51", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:35:56", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:36:01", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("18:36:06", "18:00", "Bathroom", "No Motion", "Inactive")
)

$pirStartRow = 167
$pirEndRow = $pirStartRow + $pirData.Count - 1

for ($i = 0; $i -lt $pirData.Count; $i++) {
    $r = $pirStartRow + $i
    $row = $pirData[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsPIR.Cells.Item($r, $c + 2).Value = $row[$c]
    }
}

# Column A ("Date") is the same literal string for every new row; write it
# as a text formula then flatten so it doesn't become a date serial.
$pirDateRange = $wsPIR.Range("A" + $pirStartRow + ":A" + $pirEndRow)
$pirDateRange.Formula = '="2026-01-30"'
$pirDateRange.Copy()
$pirDateRange.PasteSpecial($xlPasteValues)

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 111-119
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

$humidityData = @(
    @("18:35:09", "18:00", "Bathroom", "86.1%", "Active"),
    @("18:35:09", "18:00", "Bathroom", "85.1%", "Active"),
    @("18:35:11", "18:00", "Bathroom", "86.1%", "Active"),
    @("18:35:16", "18:00", "Bathroom", "85.1%", "Active"),
    @("18:35:22", "18:00", "Bathroom", "86.0%", "Active"),
    @("18:35:27", "18:00", "Bathroom", "86.0%", "Active"),
    @("18:35:42", "18:00", "Bathroom", "86.1%", "Active"),
    @("18:35:52", "18:00", "Bathroom", "86.0%", "Active"),
    @("18:36:02", "18:00", "Bathroom", "86.1%", "Active")
)

$humStartRow = 111
$humEndRow = $humStartRow + $humidityData.Count - 1

for ($i = 0; $i -lt $humidityData.Count; $i++) {
    $r = $humStartRow + $i
    $row = $humidityData[$i]
    $wsHumidity.Cells.Item($r, 2).Value = $row[0]
    $wsHumidity.Cells.Item($r, 3).Value = $row[1]
    $wsHumidity.Cells.Item($r, 4).Value = $row[2]
    # "Value" column (E) holds a percentage-looking string; write as a text
    # formula for now and flatten below along with the date column.
    $wsHumidity.Cells.Item($r, 5).Formula = '="' + $row[3] + '"'
    $wsHumidity.Cells.Item($r, 6).Value = $row[4]
}

# Column A ("Date") - same literal string for every new row.
$humDateRange = $wsHumidity.Range("A" + $humStartRow + ":A" + $humEndRow)
$humDateRange.Formula = '="2026-01-30"'
$humDateRange.Copy()
$humDateRange.PasteSpecial($xlPasteValues)

# Column E ("Value") - flatten the per-row text formulas to plain strings.
# NOTE: done as its own Copy/PasteSpecial (not unioned with the date range)
# because PasteSpecial on a multi-area Range mis-targets the second area.
$humValueRange = $wsHumidity.Range("E" + $humStartRow + ":E" + $humEndRow)
$humValueRange.Copy()
$humValueRange.PasteSpecial($xlPasteValues)
